$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2 and 3 with new product data ---
# Row 2
$ws.Range("B2").Value = 'Onvo OV50F900 Frameless 4K Ultra HD 50" 127 Ekran Uydu Alıcılı Android Smart LED TV'
$ws.Range("C2").Value = 'Trendyol/Beyaz Live '
$ws.Range("D2").Value = '5.649,00 TL'
$ws.Range("E2").Value = 'Trendyol/Kargomarket '
$ws.Range("F2").Value = '5.689,00 TL'
$ws.Range("H2").Value = 'https://www.akakce.com/televizyon/en-ucuz-onvo-onvo-ov50f900-127-cm-frameless-4k-uhd-lisansli-android-smart-bluetooth-dahili-uydu-alicili-1000-hz-fiyati,1623302352.html'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '0.70'
$ws.Range("G2").Style = "Normal"

# Row 3
$ws.Range("B3").Value = 'Huawei MateBook D16 i5-12450H 8 GB 512 GB SSD UHD Graphics 16" Notebook'
$ws.Range("C3").Value = 'Trendyol/HIZLIALTEKNOLOJI '
$ws.Range("D3").Value = '15.599,00 TL'
$ws.Range("E3").Value = 'Trendyol/Central Teknoloji '
$ws.Range("F3").Value = '15.629,00 TL'
$ws.Range("H3").Value = 'https://www.akakce.com/laptop-notebook/en-ucuz-huawei-matebook-d16-i5-12450h-8-gb-512-gb-ssd-uhd-graphics-16-notebook-fiyati,1954151588.html'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '0.19'
$ws.Range("G3").Style = "Normal"

# --- Append new rows 4 through 8 ---
# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 'Sapphire NITRO+ AMD RX 6700 XT 11306-01-20G 192 Bit GDDR6 12 GB Ekran Kartı'
$ws.Range("C4").Value = 'Sinerji Bilgisayar '
$ws.Range("D4").Value = '8.899,17 TL'
$ws.Range("E4").Value = 'Trendyol/Dali Teknoloji '
$ws.Range("F4").Value = '10.499,00 TL'
$ws.Range("H4").Value = 'https://www.akakce.com/ekran-karti/en-ucuz-sapphire-nitro-amd-rx-6700-xt-11306-01-20g-192-bit-gddr6-12-gb-fiyati,1100803768.html'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '15.24'
$ws.Range("G4").Style = "Normal"

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 'Lenovo Ideapad 3 82H802RKTX i3-1115G4 8 GB 256 GB SSD UHD Graphıcs 15.6" Full HD Notebook'
$ws.Range("C5").Value = 'Trendyol/Teknosa '
$ws.Range("D5").Value = '7.995,00 TL'
$ws.Range("E5").Value = 'Trendyol/VATAN BİLGİSAYAR '
$ws.Range("F5").Value = '7.999,00 TL'
$ws.Range("H5").Value = 'https://www.akakce.com/laptop-notebook/en-ucuz-lenovo-ideapad-3-82h802rktx-i3-1115g4-8-gb-256-gb-ssd-uhd-graphics-15-6-full-hd-notebook-fiyati,2110900082.html'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '0.05'
$ws.Range("G5").Style = "Normal"

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 'LG 55QNED7S6QA 4K Ultra HD 55" 140 Ekran Uydu Alıcılı Smart QNED TV'
$ws.Range("C6").Value = 'Pttavm/KIRMIZI ELMA '
$ws.Range("D6").Value = '16.619,00 TL'
$ws.Range("E6").Value = 'Trendyol/Teknomix '
$ws.Range("F6").Value = '16.865,00 TL'
$ws.Range("H6").Value = 'https://www.akakce.com/televizyon/en-ucuz-lg-55qned7s6-55inc-139-cm-4k-uhd-webos-smart-tv-uydu-alicili-fiyati,85415069.html'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '1.46'
$ws.Range("G6").Style = "Normal"

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 'Daikin MC70L Hava Temizleme Cihazı'
$ws.Range("C7").Value = 'N11/basaranstore '
$ws.Range("D7").Value = '3.707,95 TL'
$ws.Range("E7").Value = 'Veyisoglugrup.com 9,9 490 Yorum'
$ws.Range("F7").Value = '3.749,00 TL'
$ws.Range("H7").Value = 'https://www.akakce.com/hava-temizleme-cihazi/en-ucuz-daikin-mc70l-fiyati,833500.html'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '1.12'
$ws.Range("G7").Style = "Normal"

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 'Sony A7 III Body Aynasız Fotoğraf Makinesi'
$ws.Range("C8").Value = 'Hepsiburada/FOTO ÇARŞI '
$ws.Range("D8").Value = '32.148,00 TL'
$ws.Range("E8").Value = 'Klas Foto 8,7 26 Yorum'
$ws.Range("F8").Value = '32.499,00 TL'
$ws.Range("H8").Value = 'https://www.akakce.com/fotograf-makinesi/en-ucuz-sony-a7-iii-body-fiyati,209457498.html'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '1.08'
$ws.Range("G8").Style = "Normal"

# --- Copy style from A2 (bold, bordered, centered header-like style) to new A4:A8 cells ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4:A8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

